$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D(4), J(10), K(11), L(12), M(13), P(16)
# Format: row -> @(D, J, K, L, M, P)
$rowData = @{
    2  = @(45091, 40, 20000, 22000, 21000, 1400)
    4  = @(45133, 50, 22000, 22000, 22000, 1467)
    5  = @(44750, 140, 19000, 20000, 19571, 1305)
    6  = @(45084, 90, 22000, 23000, 22556, 1504)
    7  = @(45119, 50, 20000, 20000, 20000, 1333)
    8  = @(45141, 50, 8500, 9000, 8800, 587)
    9  = @(44749, 90, 17000, 18000, 17556, 1170)
    10 = @(45063, 40, 21000, 22000, 21500, 1433)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
    $ws.Cells.Item($row, 11).Value = $vals[2]
    $ws.Cells.Item($row, 12).Value = $vals[3]
    $ws.Cells.Item($row, 13).Value = $vals[4]
    $ws.Cells.Item($row, 16).Value = $vals[5]
}
